$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 / J1 - copy style (bold/centered/border) from existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for column I (I0) and column J (IF), rows 2-27
$values = @(
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(10, 10),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(5, 5),
    @(4, 4),
    @(9, 9),
    @(9, 9),
    @(5, 7),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(3, 3),
    @(6, 6),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
